$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Морозильники")
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""
